# Refactor lecturer and module timetable without redundant version
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the now-unused sample/test values from the availability grid
$ws.Range("C3").ClearContents()
$ws.Range("A4").ClearContents()
$ws.Range("E4").ClearContents()

# Move the active selection to C7 (matches saved sheet view state)
$ws.Range("C7").Select()
